$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("service")

# Add three new indicator rows (Stunting) below the existing data (rows 127-129)
$ws.Range("A127").Value = "IHlMRYVhmOX"
$ws.Range("B127").Value = "MOH 711 Stunting 0-<6 months"
$ws.Range("C127").Value = "<5yrs who are stunted"
$ws.Range("D127").Value = "stunted"
$ws.Range("E127").Value = 0.25

$ws.Range("A128").Value = "madOGCuPg8q"
$ws.Range("B128").Value = "MOH 711 Stunting 24-59 Months"
$ws.Range("C128").Value = "<5yrs who are stunted"
$ws.Range("D128").Value = "stunted"
$ws.Range("E128").Value = 0.25

$ws.Range("A129").Value = "f0SJMmqh8un"
$ws.Range("B129").Value = "MOH 711 Stunting 6-23 months"
$ws.Range("C129").Value = "<5yrs who are stunted"
$ws.Range("D129").Value = "stunted"
$ws.Range("E129").Value = 0.25

# Selection state to match the saved view position after editing near the bottom of the sheet
$ws.Range("C132").Select()
